# Edit generated for Figure_3_SourceData.xlsx
# Rewrites sheet '3A' data table (new columns/rows) and adds column G 'Plates' to sheet '3B'.

$wb = $excel.ActiveWorkbook

# ---------- Sheet "3A": full rebuild of the table ----------
$wsA = $wb.Worksheets.Item("3A")

# Clear out the old table body (rows below the new 21-row extent) so no stale cells remain.
$wsA.Range("A22:D29").ClearContents()

# Copy the existing bold/bordered header style (A1) onto the two new header columns (E1:F1)
# before writing their captions, so the whole header row keeps a uniform look.
$wsA.Range("A1").Copy()
$wsA.Range("E1:F1").PasteSpecial(-4122)
$wsA.Application.CutCopyMode = $false

$headerA = @("Plate Name", "Clone type", "Density", "Cell Type", "Percent Replicating", "Percent Matching")
for ($c = 1; $c -le $headerA.Length; $c++) {
    $wsA.Cells.Item(1, $c).Value = $headerA[$c - 1]
}

$dataA = @(
    @('BR00117105', 'Parental', 'Lower', 'A549', 62.22222222222222, 18.6046511627907),
    @('BR00117106', 'Parental', 'Lower', 'A549', 67.77777777777779, 16.27906976744186),
    @('BR00117107', 'Parental', 'Lower', 'A549', 60, 18.6046511627907),
    @('BR00117108', 'Parental', 'Higher', 'A549', 58.88888888888889, 16.27906976744186),
    @('BR00117097', 'Polyclonal', 'Lower', 'A549', 60, 20.93023255813954),
    @('BR00117098', 'Polyclonal', 'Higher', 'A549', 64.44444444444444, 18.6046511627907),
    @('BR00117100', 'Monoclonal 1', 'Lower', 'A549', 64.44444444444444, 20.93023255813954),
    @('BR00117099', 'Monoclonal 1', 'Higher', 'A549', 57.77777777777777, 13.95348837209302),
    @('BR00117104', 'Monoclonal 2', 'Lower', 'A549', 68.88888888888889, 16.27906976744186),
    @('BR00117101', 'Monoclonal 2', 'Higher', 'A549', 68.88888888888889, 16.27906976744186),
    @('BR00117110', 'Monoclonal 3', 'Lower', 'A549', 55.55555555555556, 13.95348837209302),
    @('BR00117109', 'Monoclonal 3', 'Higher', 'A549', 62.22222222222222, 18.6046511627907),
    @('BR00117093', 'Parental', 'Lower', 'U2OS', 65.55555555555556, 25.58139534883721),
    @('BR00117094', 'Parental', 'Lower', 'U2OS', 64.04494382022472, 23.80952380952381),
    @('BR00117095', 'Parental', 'Lower', 'U2OS', 65.55555555555556, 20.93023255813954),
    @('BR00117096', 'Parental', 'Higher', 'U2OS', 64.44444444444444, 18.6046511627907),
    @('BR00117088', 'Polyclonal', 'Lower', 'U2OS', 59.55056179775281, 23.80952380952381),
    @('BR00117089', 'Polyclonal', 'Higher', 'U2OS', 61.79775280898876, 16.66666666666666),
    @('BR00117091', 'Monoclonal 1', 'Lower', 'U2OS', 63.33333333333333, 18.6046511627907),
    @('BR00117092', 'Monoclonal 1', 'Higher', 'U2OS', 61.11111111111112, 16.27906976744186)
)

for ($r = 0; $r -lt $dataA.Length; $r++) {
    $row = $dataA[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsA.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ---------- Sheet "3B": append the "Plates" column (G) ----------
$wsB = $wb.Worksheets.Item("3B")

# Copy the header style from an existing header cell (A1) onto the new G1 header cell.
$wsB.Range("A1").Copy()
$wsB.Range("G1").PasteSpecial(-4122)
$wsB.Application.CutCopyMode = $false
$wsB.Cells.Item(1, 7).Value = "Plates"

$platesB = @(
    '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'']',
    '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00117020'', ''BR00117021'']',
    '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00118050'', ''BR00117006'']',
    '[''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'', ''BR00117020'', ''BR00117021'']',
    '[''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'', ''BR00118050'', ''BR00117006'']',
    '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'']',
    '[''BR00116991'', ''BR00116992'', ''BR00116993'', ''BR00116994'', ''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'']',
    '[''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'', ''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'']',
    '[''BR00117015'', ''BR00117016'', ''BR00117017'', ''BR00117019'', ''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'']',
    '[''BR00117020'', ''BR00117021'', ''BR00118050'', ''BR00117006'']',
    '[''BR00117020'', ''BR00117021'', ''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'']',
    '[''BR00117020'', ''BR00117021'', ''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'']',
    '[''BR00118050'', ''BR00117006'', ''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'']',
    '[''BR00118050'', ''BR00117006'', ''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'']',
    '[''BR00118041'', ''BR00118042'', ''BR00118043'', ''BR00118044'', ''BR00117000'', ''BR00117003'', ''BR00117004'', ''BR00117005'']',
    '[''BR00116995'', ''BR00117024'', ''BR00117025'', ''BR00117026'', ''BR00117010'', ''BR00117011'', ''BR00117012'', ''BR00117013'']',
    '[''BR00116995'', ''BR00117024'', ''BR00117025'', ''BR00117026'', ''BR00117022'', ''BR00117023'']',
    '[''BR00116995'', ''BR00117024'', ''BR00117025'', ''BR00117026'', ''BR00118039'', ''BR00118040'']',
    '[''BR00117010'', ''BR00117011'', ''BR00117012'', ''BR00117013'', ''BR00117022'', ''BR00117023'']',
    '[''BR00117010'', ''BR00117011'', ''BR00117012'', ''BR00117013'', ''BR00118039'', ''BR00118040'']',
    '[''BR00116995'', ''BR00117024'', ''BR00117025'', ''BR00117026'', ''BR00118045'', ''BR00118046'', ''BR00118047'', ''BR00118048'']',
    '[''BR00116995'', ''BR00117024'', ''BR00117025'', ''BR00117026'', ''BR00116996'', ''BR00116997'', ''BR00116998'', ''BR00116999'']',
    '[''BR00117010'', ''BR00117011'', ''BR00117012'', ''BR00117013'', ''BR00118045'', ''BR00118046'', ''BR00118047'', ''BR00118048'']',
    '[''BR00117010'', ''BR00117011'', ''BR00117012'', ''BR00117013'', ''BR00116996'', ''BR00116997'', ''BR00116998'', ''BR00116999'']',
    '[''BR00117022'', ''BR00117023'', ''BR00118039'', ''BR00118040'']',
    '[''BR00117022'', ''BR00117023'', ''BR00118045'', ''BR00118046'', ''BR00118047'', ''BR00118048'']',
    '[''BR00117022'', ''BR00117023'', ''BR00116996'', ''BR00116997'', ''BR00116998'', ''BR00116999'']',
    '[''BR00118039'', ''BR00118040'', ''BR00118045'', ''BR00118046'', ''BR00118047'', ''BR00118048'']',
    '[''BR00118039'', ''BR00118040'', ''BR00116996'', ''BR00116997'', ''BR00116998'', ''BR00116999'']',
    '[''BR00118045'', ''BR00118046'', ''BR00118047'', ''BR00118048'', ''BR00116996'', ''BR00116997'', ''BR00116998'', ''BR00116999'']'
)

for ($r = 0; $r -lt $platesB.Length; $r++) {
    $wsB.Cells.Item($r + 2, 7).Value = $platesB[$r]
}

Write-Output "3A rows written: $($dataA.Length); 3B Plates rows written: $($platesB.Length)"
